$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F9").Value = 628
$ws.Range("F11").Value = 468
$ws.Range("F12").Value = 749
$ws.Range("F14").Value = 1214
$ws.Range("F17").Value = 1237
$ws.Range("F18").Value = 298
$ws.Range("F19").Value = 1585
$ws.Range("F21").Value = 726
$ws.Range("F22").Value = 324
$ws.Range("F25").Value = 1372
$ws.Range("F27").Value = 76
$ws.Range("F31").Value = 258659
$ws.Range("F32").Value = 972
$ws.Range("F33").Value = 10
$ws.Range("F35").Value = 1318
$ws.Range("F36").Value = 828
$ws.Range("F42").Value = 832

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 144
$ws.Range("F12").Value = 70
$ws.Range("F14").Value = 2530
$ws.Range("F21").Value = 17
$ws.Range("F23").Value = 434
$ws.Range("F34").Value = 132

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2760
$ws.Range("F6").Value = 4519
$ws.Range("F11").Value = 422
$ws.Range("F12").Value = 210
$ws.Range("F15").Value = 443

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2760
$ws.Range("F6").Value = 4519
$ws.Range("F9").Value = 210
$ws.Range("F10").Value = 210
$ws.Range("F17").Value = 628
$ws.Range("F18").Value = 468
$ws.Range("F19").Value = 749
$ws.Range("F20").Value = 2530
$ws.Range("F23").Value = 1214
$ws.Range("F25").Value = 1237
$ws.Range("F27").Value = 298
$ws.Range("F29").Value = 1585
$ws.Range("F31").Value = 324
$ws.Range("F32").Value = 443
$ws.Range("F33").Value = 434
$ws.Range("F34").Value = 1372
$ws.Range("F40").Value = 972
$ws.Range("F43").Value = 832

